$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 8682.200000000001
$ws.Range("I34").Value = 9352.75
$ws.Range("J34").Value = 6000
$ws.Range("K34").Value = 9352.75
$ws.Range("L34").Value = 6000
$ws.Range("M34").Value = -9149.75
$ws.Range("N34").Value = -6406

$ws.Range("H36").Value = 8682.200000000001
$ws.Range("I36").Value = 9352.75
$ws.Range("J36").Value = 6000
$ws.Range("K36").Value = 9352.75
$ws.Range("L36").Value = 6000
$ws.Range("M36").Value = -8637.75
$ws.Range("N36").Value = -7430

$ws.Range("H38").Value = 414.83334

$ws.Range("H43").Value = 1600.0938
$ws.Range("I43").Value = 1419.3462
$ws.Range("J43").Value = 2383.3333
$ws.Range("K43").Value = 1419.3462
$ws.Range("L43").Value = 2383.3333
$ws.Range("M43").Value = -1350.3462
$ws.Range("N43").Value = -2521.3333

$ws.Range("H51").Value = 9813839
$ws.Range("I51").Value = 55561424
$ws.Range("J51").Value = 10785.714
$ws.Range("K51").Value = 55561424
$ws.Range("L51").Value = 10785.714
$ws.Range("M51").Value = -55560940
$ws.Range("N51").Value = -11753.714

$ws.Range("H58").Value = 1010

$ws.Range("H80").Value = 22223752
$ws.Range("I80").Value = 55556484
$ws.Range("K80").Value = 166669452
$ws.Range("M80").Value = -166668454

$ws.Range("H83").Value = 22223752
$ws.Range("I83").Value = 55556484
$ws.Range("K83").Value = 500008356
$ws.Range("M83").Value = -500003364

$ws.Range("H86").Value = 2115.6897
$ws.Range("I86").Value = 2240.1177
$ws.Range("J86").Value = 1939.4166
$ws.Range("K86").Value = 2240.1177
$ws.Range("L86").Value = 1939.4166
$ws.Range("M86").Value = -1117.1177
$ws.Range("N86").Value = -4185.4166

$ws.Range("H88").Value = 528598.6
$ws.Range("I88").Value = 1948444.2
$ws.Range("J88").Value = 12291.091
$ws.Range("K88").Value = 1948444.2
$ws.Range("L88").Value = 12291.091
$ws.Range("M88").Value = -1948038.2
$ws.Range("N88").Value = -13103.091

$ws.Range("H89").Value = 2115.6897
$ws.Range("I89").Value = 2240.1177
$ws.Range("J89").Value = 1939.4166
$ws.Range("K89").Value = 11200.5885
$ws.Range("L89").Value = 9697.083000000001
$ws.Range("M89").Value = -5584.588499999998
$ws.Range("N89").Value = -20929.083

$ws.Range("H91").Value = 528598.6
$ws.Range("I91").Value = 1948444.2
$ws.Range("J91").Value = 12291.091
$ws.Range("K91").Value = 1948444.2
$ws.Range("L91").Value = 12291.091
$ws.Range("M91").Value = -1947040.2
$ws.Range("N91").Value = -15099.091

$ws.Range("H113").Value = 2399.75
$ws.Range("I113").Value = 2266.3333
$ws.Range("J113").Value = 2800
$ws.Range("K113").Value = 2266.3333
$ws.Range("L113").Value = 2800
$ws.Range("M113").Value = 987.6667000000002
$ws.Range("N113").Value = -9308

$ws.Range("H125").Value = 1792.5714
$ws.Range("I125").Value = 1554.6666
$ws.Range("J125").Value = 1971
$ws.Range("K125").Value = 13991.9994
$ws.Range("L125").Value = 17739
$ws.Range("M125").Value = -11531.9994
$ws.Range("N125").Value = -22659

$ws.Range("H141").Value = 3221.913
$ws.Range("I141").Value = 1840
$ws.Range("J141").Value = 4729.4546
$ws.Range("K141").Value = 5520
$ws.Range("L141").Value = 14188.3638
$ws.Range("M141").Value = -340
$ws.Range("N141").Value = -24548.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1501.1887
$ws.Range("I61").Value = 984.25
$ws.Range("J61").Value = 3091.7693
$ws.Range("K61").Value = 984.25
$ws.Range("L61").Value = 3091.7693
$ws.Range("M61").Value = -772.25
$ws.Range("N61").Value = -3515.7693

$ws.Range("H122").Value = 1868.2222
$ws.Range("I122").Value = 1537.6364
$ws.Range("J122").Value = 2387.7144
$ws.Range("K122").Value = 4612.9092
$ws.Range("L122").Value = 7163.1432
$ws.Range("M122").Value = -2162.9092
$ws.Range("N122").Value = -12063.1432

$ws.Range("H136").Value = 1501.1887
$ws.Range("I136").Value = 984.25
$ws.Range("J136").Value = 3091.7693
$ws.Range("K136").Value = 2952.75
$ws.Range("L136").Value = 9275.3079
$ws.Range("M136").Value = -402.75
$ws.Range("N136").Value = -14375.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7777
$ws.Range("I86").Value = 3333
$ws.Range("K86").Value = 3333
$ws.Range("M86").Value = -2210

$ws.Range("H89").Value = 7777
$ws.Range("I89").Value = 3333
$ws.Range("K89").Value = 16665
$ws.Range("M89").Value = -11049

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3972042.8
$ws.Range("I31").Value = 1812
$ws.Range("J31").Value = 6807921.5
$ws.Range("K31").Value = 1812
$ws.Range("L31").Value = 6807921.5
$ws.Range("M31").Value = -1517
$ws.Range("N31").Value = -6808511.5

$ws.Range("H34").Value = 3972042.8
$ws.Range("I34").Value = 1812
$ws.Range("J34").Value = 6807921.5
$ws.Range("K34").Value = 1812
$ws.Range("L34").Value = 6807921.5
$ws.Range("M34").Value = -1610
$ws.Range("N34").Value = -6808325.5

$ws.Range("H58").Value = 2420.4482
$ws.Range("I58").Value = 1506.4286
$ws.Range("J58").Value = 3273.5334
$ws.Range("K58").Value = 1506.4286
$ws.Range("L58").Value = 3273.5334
$ws.Range("M58").Value = -1303.4286
$ws.Range("N58").Value = -3679.5334

$ws.Range("H136").Value = 2420.4482
$ws.Range("I136").Value = 1506.4286
$ws.Range("J136").Value = 3273.5334
$ws.Range("K136").Value = 4519.2858
$ws.Range("L136").Value = 9820.600199999999
$ws.Range("M136").Value = -1969.2858
$ws.Range("N136").Value = -14920.6002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1809.5
$ws.Range("I59").Value = 1149.375
$ws.Range("J59").Value = 4450
$ws.Range("K59").Value = 3448.125
$ws.Range("L59").Value = 13350
$ws.Range("M59").Value = -2908.125
$ws.Range("N59").Value = -14430

$ws.Range("H120").Value = 384374.88
$ws.Range("J120").Value = 11999.833
$ws.Range("L120").Value = 35999.499
$ws.Range("N120").Value = -45675.499

$ws.Range("H126").Value = 4164.7144
$ws.Range("J126").Value = 4164.7144
$ws.Range("L126").Value = 12494.1432
$ws.Range("N126").Value = -22374.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7667
$ws.Range("I80").Value = 5889.4443
$ws.Range("J80").Value = 10333.333
$ws.Range("K80").Value = 5889.4443
$ws.Range("L80").Value = 10333.333
$ws.Range("M80").Value = -4891.4443
$ws.Range("N80").Value = -12329.333

$ws.Range("H83").Value = 7667
$ws.Range("I83").Value = 5889.4443
$ws.Range("J83").Value = 10333.333
$ws.Range("K83").Value = 29447.2215
$ws.Range("L83").Value = 51666.665
$ws.Range("M83").Value = -24455.2215
$ws.Range("N83").Value = -61650.665

$ws.Range("H102").Value = 1206
$ws.Range("I102").Value = 941.3333
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 941.3333
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 680.6667
$ws.Range("N102").Value = -5244

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3512.2144
$ws.Range("I16").Value = 3333.7273
$ws.Range("J16").Value = 4166.6665
$ws.Range("K16").Value = 3333.7273
$ws.Range("L16").Value = 4166.6665
$ws.Range("M16").Value = -3163.7273
$ws.Range("N16").Value = -4506.6665

$ws.Range("H68").Value = 3673.5557
$ws.Range("J68").Value = 3713.0667
$ws.Range("L68").Value = 3713.0667
$ws.Range("N68").Value = -5211.066699999999

$ws.Range("H71").Value = 3673.5557
$ws.Range("J71").Value = 3713.0667
$ws.Range("L71").Value = 18565.3335
$ws.Range("N71").Value = -26053.3335

$ws.Range("H82").Value = 8334168.5
$ws.Range("I82").Value = 986
$ws.Range("J82").Value = 11905533
$ws.Range("K82").Value = 986
$ws.Range("L82").Value = 11905533
$ws.Range("M82").Value = -625
$ws.Range("N82").Value = -11906255

$ws.Range("H85").Value = 8334168.5
$ws.Range("I85").Value = 986
$ws.Range("J85").Value = 11905533
$ws.Range("K85").Value = 986
$ws.Range("L85").Value = 11905533
$ws.Range("M85").Value = 262
$ws.Range("N85").Value = -11908029

$ws.Range("H122").Value = 79300.30499999999
$ws.Range("J122").Value = 3250
$ws.Range("L122").Value = 9750
$ws.Range("N122").Value = -14650

$ws.Range("H132").Value = 3147.5107
$ws.Range("I132").Value = 2451.5557
$ws.Range("J132").Value = 4087.05
$ws.Range("K132").Value = 7354.6671
$ws.Range("L132").Value = 12261.15
$ws.Range("M132").Value = -4824.6671
$ws.Range("N132").Value = -17321.15

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 833.6667
$ws.Range("I81").Value = 600
$ws.Range("J81").Value = 1301
$ws.Range("K81").Value = 1200
$ws.Range("L81").Value = 2602
$ws.Range("M81").Value = -139
$ws.Range("N81").Value = -4724

$ws.Range("H84").Value = 833.6667
$ws.Range("I84").Value = 600
$ws.Range("J84").Value = 1301
$ws.Range("K84").Value = 6000
$ws.Range("L84").Value = 13010
$ws.Range("M84").Value = -696
$ws.Range("N84").Value = -23618

$ws.Range("H107").Value = 5264150.5
$ws.Range("I107").Value = 918.7143
$ws.Range("J107").Value = 20001200
$ws.Range("K107").Value = 2756.1429
$ws.Range("L107").Value = 60003600
$ws.Range("M107").Value = -836.1428999999998
$ws.Range("N107").Value = -60007440

$ws.Range("H123").Value = 39714.5
$ws.Range("J123").Value = 39714.5
$ws.Range("L123").Value = 39714.5
$ws.Range("N123").Value = -49514.5

$ws.Range("H136").Value = 389690.22
$ws.Range("I136").Value = 543254.1
$ws.Range("J136").Value = 1263.8235
$ws.Range("K136").Value = 1629762.3
$ws.Range("L136").Value = 3791.4705
$ws.Range("M136").Value = -1627212.3
$ws.Range("N136").Value = -8891.470499999999
